# feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" worksheet (fund-holdings detail, same shape as the
# 2021-Q3 / 2021-Q4 sheets) and inserts a corresponding summary row at the
# top of the "总计" (totals) sheet, pushing the existing 2021-Q4 / 2021-Q3
# rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a Range as *text*, never letting Excel's
# "smart" type detection turn a numeric-looking string (e.g. a fund code
# like "007130", or a percentage-like "40.99") into a real number. A
# leading apostrophe forces text entry; resetting the Style back to
# "Normal" afterwards drops the "Number Stored as Text" quote-prefix
# flag Excel would otherwise leave behind, so no stray cell style is
# introduced.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# Helper: copy the formatting (only) of one cell/range onto another so
# new cells can pick up the workbook's existing bold/border "header /
# index column" look (style index 2) without fabricating a brand new
# style entry.
# ---------------------------------------------------------------------
function Copy-Format {
    param($fromRange, $toRange)
    $fromRange.Copy()
    $toRange.PasteSpecial(-4122)   # xlPasteFormats
}

# =======================================================================
# 1. Duplicate the "总计" sheet and place the copy *after* it, then swap
#    names so that:
#      - the sheet now in "总计"'s old slot becomes "2022-Q1"
#      - the new copy (last position) becomes the new "总计"
#    This reproduces the exact sheetId/r:id pattern of the authored
#    change (2022-Q1 keeps the original sheetId/rId, 总计 gets new ones).
# =======================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Copy($null, $totalSheet)

$fundSheet = $wb.Worksheets.Item(3)
$newTotalSheet = $wb.Worksheets.Item(4)

$fundSheet.Name = "2022-Q1"
$newTotalSheet.Name = "总计"

# =======================================================================
# 2. Build the "2022-Q1" fund-holdings sheet from scratch (same layout
#    as the 2021-Q3 / 2021-Q4 sheets: 基金代码/基金名称/基金规模/
#    股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# =======================================================================
$fundSheet.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $fundSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$fundRows = @(
    @("007130", "中庚小盘价值股票",   "40.99", "93.10", "2.54", "1.0411", 10),
    @("090010", "大成中证红利指数A",  "34.51", "93.73", "1.59", "0.5487", 6),
    @("515180", "易方达中证红利ETF",  "16.55", "99.58", "1.69", "0.2797", 6),
    @("515080", "招商中证红利ETF",    "9.06",  "99.25", "1.68", "0.1522", 6),
    @("007801", "大成中证红利指数C",  "3.87",  "93.73", "1.59", "0.0615", 6),
    @("515890", "博时中证红利ETF",    "2.59",  "98.55", "1.67", "0.0433", 6),
    @("161907", "万家中证红利指数(LOF)", "1.34", "94.87", "1.61", "0.0216", 7)
)

$row = 2
foreach ($fund in $fundRows) {
    $fundSheet.Cells.Item($row, 1).Value = $row - 2
    Set-TextValue $fundSheet.Cells.Item($row, 2) $fund[0]
    $fundSheet.Cells.Item($row, 3).Value = $fund[1]
    Set-TextValue $fundSheet.Cells.Item($row, 4) $fund[2]
    Set-TextValue $fundSheet.Cells.Item($row, 5) $fund[3]
    Set-TextValue $fundSheet.Cells.Item($row, 6) $fund[4]
    Set-TextValue $fundSheet.Cells.Item($row, 7) $fund[5]
    $fundSheet.Cells.Item($row, 8).Value = $fund[6]
    $row++
}

# Header row (B1:H1) and index column (A2:A8) reuse the bold/border
# "s=2" style already present elsewhere in the workbook (e.g. 2021-Q4).
Copy-Format $wb.Worksheets.Item("2021-Q4").Range("B1") $fundSheet.Range("B1:H1")
Copy-Format $wb.Worksheets.Item("2021-Q4").Range("A2") $fundSheet.Range("A2:A8")

# =======================================================================
# 3. Insert the 2022-Q1 summary as the new first data row of "总计",
#    shifting 2021-Q4 / 2021-Q3 down by one row.
# =======================================================================
$newTotalSheet.Range("B4").Value = $newTotalSheet.Range("B3").Value2
$newTotalSheet.Range("C4").Value = $newTotalSheet.Range("C3").Value2
$newTotalSheet.Range("D4").Value = $newTotalSheet.Range("D3").Value2
$newTotalSheet.Range("A4").Value = 2

$newTotalSheet.Range("B3").Value = $newTotalSheet.Range("B2").Value2
$newTotalSheet.Range("C3").Value = $newTotalSheet.Range("C2").Value2
$newTotalSheet.Range("D3").Value = $newTotalSheet.Range("D2").Value2
$newTotalSheet.Range("A3").Value = 1

$newTotalSheet.Range("A2").Value = 0
$newTotalSheet.Range("B2").Value = "2022-Q1"
$newTotalSheet.Range("C2").Value = 7
$newTotalSheet.Range("D2").Value = 2.15

# A4 is a brand-new cell (the sheet used to stop at row 3), so it needs
# the index-column style copied onto it explicitly.
Copy-Format $newTotalSheet.Range("A2") $newTotalSheet.Range("A4")

# =======================================================================
# 4. Restore the original active sheet/tab (2021-Q3), which the sheet
#    Copy()/rename operations above would otherwise have moved.
# =======================================================================
$wb.Worksheets.Item("2021-Q3").Activate()
$wb.Worksheets.Item("2021-Q3").Range("A1").Select()
